# Updated symbol list on Wed Dec 14 12:59:50 UTC 2022 with GitHub Actions
#
# This script reproduces the data refresh applied to the cryptocurrency
# price sheet: most rows simply get a refreshed "Price" (column D) value,
# while rows 15-26 are a block that shifted down by one because
# "ProBitToken" moved up from row 26 to row 15 (pushing the other rows
# down) - each of those rows also carries its own refreshed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows whose Price (column D) value simply gets refreshed in place.
# ---------------------------------------------------------------------
$priceUpdates = @{
    "D2"  = "271.69"
    "D3"  = "22.96"
    "D4"  = "6.375"
    "D5"  = "0.06227"
    "D6"  = "3.641"
    "D7"  = "6.702"
    "D8"  = "1.372"
    "D9"  = "0.8356"
    "D10" = "0.01377"
    "D11" = "0.1633"
    "D12" = "0.08366"
    "D13" = "0.03409"
    "D14" = "0.03103"
    "D40" = "0.04681"
    "D41" = "0.006958"
    "D43" = "0.003425"
    "D44" = "0.01271"
    "D45" = "0.00006255"
    "D47" = "0.8995"
    "D48" = "0.06055"
}

foreach ($addr in $priceUpdates.Keys) {
    # Force the cell to keep a text type (not get reinterpreted as a
    # number) so that values such as "0.006210" or "0.00006255" keep
    # their exact textual representation (trailing zeros, no scientific
    # notation, etc.), then restore the default "Normal" style so the
    # cell doesn't end up with a stray number-format style applied.
    $range = $ws.Range($addr)
    $range.NumberFormat = "@"
    $range.Value = $priceUpdates[$addr]
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Rows 15-26: full row content refresh (Coin / Link / Price / Volume).
# ProBitToken jumps from row 26 up to row 15, shifting BitMartToken,
# MCDex, BitForexToken, CoinExToken, TigerCash, BitKan, HotbitToken,
# NitroEx, LEO, BTSEToken and BitpandaEcosystemToken down by one row
# each, with their own data refreshed as well.
# ---------------------------------------------------------------------
$rows = @(
    @{ Row = 15; Coin = "ProBitToken";              Link = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob";                Price = "0.1271";    Volume = "14ProBitTokenPROB" }
    @{ Row = 16; Coin = "BitMartToken";              Link = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx";            Price = "0.09325";   Volume = "15BitMartTokenBMX" }
    @{ Row = 17; Coin = "MCDex";                     Link = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb";                       Price = "3.885";     Volume = "16MCDexMCB" }
    @{ Row = 18; Coin = "BitForexToken";             Link = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf";            Price = "0.001651";  Volume = "17BitForexTokenBF" }
    @{ Row = 19; Coin = "CoinExToken";                Link = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet";            Price = "0.04813";   Volume = "18CoinExTokenCET" }
    @{ Row = 20; Coin = "TigerCash";                 Link = "https://coinranking.com/coin/6hIn06L2+tigercash-tch";                    Price = "0.006210";  Volume = "19TigerCashTCH" }
    @{ Row = 21; Coin = "BitKan";                    Link = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan";                  Price = "0.001087";  Volume = "20BitKanKAN" }
    @{ Row = 22; Coin = "HotbitToken";                Link = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb";             Price = "0.003331";  Volume = "21HotbitTokenHTBWorstin24h" }
    @{ Row = 23; Coin = "NitroEx";                   Link = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx";                   Price = "0.0001499"; Volume = "22NitroExNTX" }
    @{ Row = 24; Coin = "LEO";                       Link = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo";                      Price = "3.731";     Volume = "23LEOLEO" }
    @{ Row = 25; Coin = "BTSEToken";                 Link = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse";               Price = "2.363";     Volume = "24BTSETokenBTSE" }
    @{ Row = 26; Coin = "BitpandaEcosystemToken";    Link = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best";      Price = "0.3401";    Volume = "25BitpandaEcosystemTokenBEST" }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.Coin
    $ws.Range("C$r").Value = $item.Link

    $dRange = $ws.Range("D$r")
    $dRange.NumberFormat = "@"
    $dRange.Value = $item.Price
    $dRange.Style = "Normal"

    $ws.Range("E$r").Value = $item.Volume
}
